$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$nl = [char]11

$tbl.Cell(1, 1).Range.Text = "63 x 79" + $nl + "  7    9" + $nl + "  ----" + $nl + "6|    |" + $nl + "3|    |"
$tbl.Cell(1, 2).Range.Text = "16 x 31" + $nl + "  3    1" + $nl + "  ----" + $nl + "1|    |" + $nl + "6|    |"
$tbl.Cell(1, 3).Range.Text = "41 x 70" + $nl + "  7    0" + $nl + "  ----" + $nl + "4|    |" + $nl + "1|    |"
$tbl.Cell(2, 1).Range.Text = "17 x 83" + $nl + "  8    3" + $nl + "  ----" + $nl + "1|    |" + $nl + "7|    |"
$tbl.Cell(2, 2).Range.Text = "29 x 34" + $nl + "  3    4" + $nl + "  ----" + $nl + "2|    |" + $nl + "9|    |"
$tbl.Cell(2, 3).Range.Text = "24 x 50" + $nl + "  5    0" + $nl + "  ----" + $nl + "2|    |" + $nl + "4|    |"
$tbl.Cell(3, 1).Range.Text = "99 x 57" + $nl + "  5    7" + $nl + "  ----" + $nl + "9|    |" + $nl + "9|    |"
$tbl.Cell(3, 2).Range.Text = "81 x 46" + $nl + "  4    6" + $nl + "  ----" + $nl + "8|    |" + $nl + "1|    |"
$tbl.Cell(3, 3).Range.Text = "69 x 96" + $nl + "  9    6" + $nl + "  ----" + $nl + "6|    |" + $nl + "9|    |"
$tbl.Cell(4, 1).Range.Text = "21 x 92" + $nl + "  9    2" + $nl + "  ----" + $nl + "2|    |" + $nl + "1|    |"
$tbl.Cell(4, 2).Range.Text = "13 x 95" + $nl + "  9    5" + $nl + "  ----" + $nl + "1|    |" + $nl + "3|    |"
$tbl.Cell(4, 3).Range.Text = "71 x 58" + $nl + "  5    8" + $nl + "  ----" + $nl + "7|    |" + $nl + "1|    |"
$tbl.Cell(5, 1).Range.Text = "65 x 79" + $nl + "  7    9" + $nl + "  ----" + $nl + "6|    |" + $nl + "5|    |"
$tbl.Cell(5, 2).Range.Text = "29 x 47" + $nl + "  4    7" + $nl + "  ----" + $nl + "2|    |" + $nl + "9|    |"
$tbl.Cell(5, 3).Range.Text = "35 x 97" + $nl + "  9    7" + $nl + "  ----" + $nl + "3|    |" + $nl + "5|    |"
